# Updates cryptos list on Mon Jul  3 09:59:25 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for the coin rows
# (2-51) with the latest scraped snapshot. Both columns hold plain text in
# the workbook (not numbers/percentages), so every write below goes through
# Set-CellText, which prefixes a leading apostrophe for values that look
# like a number (e.g. "1.001", "0.4833") -- exactly what typing such a
# value into Excel with quote-prefix does -- so it is stored as text and
# not silently coerced into a numeric cell. Values that already read as
# text on their own (thousand-dotted prices like "30.638.04", or the
# space-padded "  +0.47%  " percentages) are written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Addr, $Text)
    if ($Text -match '^[+-]?[0-9]+(\.[0-9]+)?([eE][+-]?[0-9]+)?$') {
        $ws.Range($Addr).Value = "'" + $Text
    } else {
        $ws.Range($Addr).Value = $Text
    }
}

# row -> [price, volume]  ($null = unchanged, leave the cell alone)
$updates = @{
    2  = @('30.638.04',    '  +0.47%  ')
    3  = @('1.958.38',     '  +2.29%  ')
    4  = @($null,          '  +0.21%  ')
    5  = @('248.10',       '  +1.34%  ')
    6  = @('1.001',        '  +0.18%  ')
    7  = @('0.4833',       $null)
    8  = @('0.2939',       '  +1.89%  ')
    9  = @('0.06776',      '  +0.86%  ')
    10 = @('110.52',       '  -0.17%  ')
    11 = @('19.42',        '  +0.99%  ')
    12 = @('1.963.16',     '  +2.53%  ')
    13 = @('0.07729',      '  +2.21%  ')
    14 = @('5.464',        '  +4.46%  ')
    15 = @('0.6869',       '  +3.14%  ')
    16 = @('292.39',       '  -3.60%  ')
    17 = @('30.653.78',    '  +0.56%  ')
    18 = @($null,          '  +2.33%  ')
    19 = @('2.219.48',     $null)
    20 = @('5.642',        '  +3.18%  ')
    21 = @('0.000007665',  '  +1.37%  ')
    22 = @('0.9999',       '  +0.18%  ')
    23 = @('1.001',        $null)
    24 = @('6.594',        '  +3.22%  ')
    25 = @('9.902',        '  +4.63%  ')
    26 = @('170.97',       '  +4.03%  ')
    27 = @('20.13',        '  -1.79%  ')
    28 = @('2.189',        '  +4.93%  ')
    29 = @('0.1072',       '  +0.15%  ')
    30 = @('1.439',        '  +2.83%  ')
    31 = @('4.686',        '  +16.69%  ')
    32 = @('4.429',        '  +6.37%  ')
    33 = @('0.05106',      '  +2.67%  ')
    34 = @('0.7758',       '  +5.85%  ')
    35 = @('1.170',        '  +2.86%  ')
    36 = @('0.02059',      '  +0.55%  ')
    37 = @('2.734',        '  +0.31%  ')
    38 = @('2.712',        '  +1.43%  ')
    39 = @('2.070',        '  +2.69%  ')
    40 = @('6.178',        '  +4.57%  ')
    41 = @('109.90',       '  -1.32%  ')
    42 = @($null,          '  +1.32%  ')
    43 = @('0.8747',       '  +1.16%  ')
    44 = @('69.92',        '  +2.04%  ')
    45 = @('1.001',        '  +0.22%  ')
    46 = @('7.388',        '  +1.56%  ')
    47 = @('0.1278',       '  +3.68%  ')
    48 = @('9.376',        '  +1.06%  ')
    49 = @($null,          '  +3.35%  ')
    50 = @('47.58',        '  -4.21%  ')
    51 = @('0.4088',       '  +2.17%  ')
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]
    if ($null -ne $price) {
        Set-CellText "D$row" $price
    }
    if ($null -ne $volume) {
        Set-CellText "E$row" $volume
    }
}
